# Apply cryptos list price/volume updates (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "38.029.44"
$ws.Range("E2").Value = "  +0.22%  "

$ws.Range("D3").Value = "2.038.34"
$ws.Range("E3").Value = "  -0.71%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").Value = "'" + "228.64"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.52%  "

$ws.Range("D6").Value = "'" + "0.607"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.59%  "

$ws.Range("D7").Value = "'" + "60.71"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.62%  "

$ws.Range("E8").Value = "  +0.06%  "

$ws.Range("E9").Value = "  -0.86%  "

$ws.Range("D10").Value = "'" + "0.0820"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.98%  "

$ws.Range("E11").Value = "  +0.59%  "

$ws.Range("B12").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C12").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D12").Value = "2.341.40"
$ws.Range("E12").Value = "  -0.60%  "

$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").Value = "'" + "14.62"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.05%  "

$ws.Range("D14").Value = "'" + "21.34"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.52%  "

$ws.Range("D15").Value = "'" + "0.766"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.68%  "

$ws.Range("D16").Value = "'" + "5.18"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.79%  "

$ws.Range("D17").Value = "2.040.00"
$ws.Range("E17").Value = "  -1.75%  "

$ws.Range("D18").Value = "37.881.15"
$ws.Range("E18").Value = "  +0.03%  "

$ws.Range("D19").Value = "'" + "69.87"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.24%  "

$ws.Range("D20").Value = "'" + "5.96"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -6.00%  "

$ws.Range("D21").Value = "0.0₃0826"
$ws.Range("E21").Value = "  -1.33%  "

$ws.Range("D22").Value = "'" + "224.78"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.13%  "

$ws.Range("D23").Value = "'" + "1.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.01%  "

$ws.Range("E24").Value = "  +0.21%  "

$ws.Range("E25").Value = "  -0.37%  "

$ws.Range("D26").Value = "'" + "9.35"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.21%  "

$ws.Range("D27").Value = "'" + "167.25"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.61%  "

$ws.Range("E28").Value = "  -1.68%  "

$ws.Range("D29").Value = "'" + "18.91"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.64%  "

$ws.Range("E30").Value = "  -3.46%  "

$ws.Range("E31").Value = "  +0.70%  "

$ws.Range("D32").Value = "'" + "2.19"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +6.33%  "

$ws.Range("D33").Value = "'" + "4.42"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.38%  "

$ws.Range("D34").Value = "'" + "4.53"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.46%  "

$ws.Range("D35").Value = "'" + "0.0607"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.27%  "

$ws.Range("D36").Value = "'" + "6.49"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +6.85%  "

$ws.Range("D37").Value = "'" + "2.28"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.00%  "

$ws.Range("D38").Value = "'" + "3.29"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.27%  "

$ws.Range("E39").Value = "  +0.09%  "

$ws.Range("E40").Value = "  +7.34%  "

$ws.Range("D41").Value = "1.527.48"
$ws.Range("E41").Value = "  +1.71%  "

$ws.Range("D42").Value = "'" + "0.0218"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.52%  "

$ws.Range("E43").Value = "  -0.75%  "

$ws.Range("D44").Value = "'" + "2.83"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.79%  "

$ws.Range("D45").Value = "'" + "0.0916"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.32%  "

$ws.Range("E46").Value = "  -2.32%  "

$ws.Range("D47").Value = "'" + "4.03"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.35%  "

$ws.Range("E48").Value = "  -0.59%  "

$ws.Range("E49").Value = "  +0.16%  "

$ws.Range("E50").Value = "  +0.36%  "

$ws.Range("D51").Value = "2.230.63"
$ws.Range("E51").Value = "  -0.60%  "
